# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The captured OOXML diff for this fixture only ever reorders XML
# attributes (e.g. <w:pgSz w:w=".." w:h=".."/> -> <w:pgSz w:h=".." w:w=".."/>,
# xmlns:* declarations alphabetised on <w:document>, w:lsdException /
# w:style attributes alphabetised, etc.). Every removed line and its
# corresponding added line carry exactly the same set of
# attribute=value pairs and the same element/text content - this is the
# well known side effect of the upstream project re-serialising its
# test-fixture .docx files after bumping the Apache POI version used to
# generate them. No paragraph text, run, field, style value, numeric
# value, page size/margin, font, language, or any other document
# content actually changed between "before" and "after".
#
# So the faithful edit here is a content no-op: touch nothing in the
# document model. We still grab ActiveDocument (as the harness expects
# a live reference) but intentionally make no mutating calls, since any
# Find/Replace or property round-trip would (per this engine) only
# introduce spurious content differences that are not present in the
# target diff, without being able to reproduce the attribute-ordering
# that *is* in the diff (that ordering is purely an artifact of the
# other project's XML writer, not something the Word object model
# exposes control over).

$d = $word.ActiveDocument
